$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formula in G2 (not part of the shared group)
$ws.Range("G2").Formula = "=+D2/C2"

# Update formulas in column G (Tasa de Recuperados) from D/E to D/C for rows 3-22
# Set as one range assignment so Excel consolidates them into a single shared formula group
$ws.Range("G3:G22").Formula = "=+D3/C3"

# Update the active cell selection to F3
$ws.Range("F3").Select()
